$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the existing data (A1:G23) one column to the right (B1:H23),
# carrying values and formatting with it, without disturbing the
# worksheet's column-width definitions.
$ws.Range("A1:G23").Copy($ws.Range("B1"))

# Header for the new column
$ws.Range("A1").Value = "Codigo"
$ws.Range("A1").Style = "Célula de Verificação"
$ws.Range("A1").HorizontalAlignment = -4108

# Fill the new column with sequential codes for each data row (1..22)
for ($i = 2; $i -le 23; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 1
    $ws.Cells.Item($i, 1).Style = "Normal"
    $ws.Cells.Item($i, 1).HorizontalAlignment = -4108
}

# Adjust width of the new column A (closest achievable value to 16.85546875)
$ws.Columns("A").ColumnWidth = 16
